$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

function Set-TextValue($cell, $value) {
    # Force the cell to Text format first so values that look like numbers,
    # dates or times (phone numbers, ids, dates, times, date-times) are stored
    # verbatim as shared-string text instead of being auto-coerced into
    # numeric / date serial values.
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Sheet1 row 2 - new schedule-notify test values.
Set-TextValue $ws1.Range("F2")  "2813215443"
Set-TextValue $ws1.Range("N2")  "2024-03-12"
Set-TextValue $ws1.Range("O2")  "02:35:55 AM"
Set-TextValue $ws1.Range("P2")  "2024-03-12 11:00:06 AM"
Set-TextValue $ws1.Range("AC2") "2024-03-12"
Set-TextValue $ws1.Range("AE2") "0952316974"
Set-TextValue $ws1.Range("AN2") "97834"
Set-TextValue $ws1.Range("AT2") "8230120282"
Set-TextValue $ws1.Range("AX2") "4803916963"

# Sheet2 row 2.
Set-TextValue $ws2.Range("F2")  "2813215443"
Set-TextValue $ws2.Range("AE2") "0952316974"
Set-TextValue $ws2.Range("AT2") "8230120282"
Set-TextValue $ws2.Range("AX2") "4803916963"

# Sheet3 row 2.
Set-TextValue $ws3.Range("F2")  "2813215443"
Set-TextValue $ws3.Range("AE2") "0952316974"
Set-TextValue $ws3.Range("AT2") "8230120282"
Set-TextValue $ws3.Range("AX2") "4803916963"

# Sheet4 row 2.
Set-TextValue $ws4.Range("F2")  "2813215443"
Set-TextValue $ws4.Range("AE2") "0952316974"
Set-TextValue $ws4.Range("AT2") "8230120282"
Set-TextValue $ws4.Range("AX2") "4803916963"
